$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 2899.5557
$ws.Range("J41").Value = 2033
$ws.Range("L41").Value = 2033
$ws.Range("N41").Value = -2913
# Row 51
$ws.Range("H51").Value = 7892.5557
$ws.Range("I51").Value = 50300
$ws.Range("J51").Value = 6261.5
$ws.Range("K51").Value = 50300
$ws.Range("L51").Value = 6261.5
$ws.Range("M51").Value = -49816
$ws.Range("N51").Value = -7229.5
# Row 141
$ws.Range("H141").Value = 7180.625
$ws.Range("I141").Value = 7241.3335
$ws.Range("K141").Value = 21724.0005
$ws.Range("M141").Value = -16544.0005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 79998.5
$ws.Range("J24").Value = 79998.5
$ws.Range("L24").Value = 79998.5
$ws.Range("N24").Value = -80746.5
# Row 61
$ws.Range("H61").Value = 16141088
$ws.Range("I61").Value = 17953138
$ws.Range("K61").Value = 17953138
$ws.Range("M61").Value = -17952926
# Row 63
$ws.Range("H63").Value = 3362.5
$ws.Range("I63").Value = 3271.4285
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 3271.4285
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -2585.4285
$ws.Range("N63").Value = -5372
# Row 66
$ws.Range("H66").Value = 3362.5
$ws.Range("I66").Value = 3271.4285
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 16357.1425
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -12925.1425
$ws.Range("N66").Value = -26864
# Row 74
$ws.Range("H74").Value = 1830.3334
$ws.Range("I74").Value = 1867.2142
$ws.Range("K74").Value = 1867.2142
$ws.Range("M74").Value = -993.2141999999999
# Row 75
$ws.Range("H75").Value = 87499.5
$ws.Range("J75").Value = 87499.5
$ws.Range("L75").Value = 87499.5
$ws.Range("N75").Value = -89247.5
# Row 77
$ws.Range("H77").Value = 1830.3334
$ws.Range("I77").Value = 1867.2142
$ws.Range("K77").Value = 9336.071
$ws.Range("M77").Value = -4968.071
# Row 78
$ws.Range("H78").Value = 87499.5
$ws.Range("J78").Value = 87499.5
$ws.Range("L78").Value = 262498.5
$ws.Range("N78").Value = -271234.5
# Row 88
$ws.Range("H88").Value = 2450.9333
$ws.Range("J88").Value = 2825.111
$ws.Range("L88").Value = 2825.111
$ws.Range("N88").Value = -3637.111
# Row 91
$ws.Range("H91").Value = 2450.9333
$ws.Range("J91").Value = 2825.111
$ws.Range("L91").Value = 2825.111
$ws.Range("N91").Value = -5633.111
# Row 97
$ws.Range("H97").Value = 4689.6665
$ws.Range("J97").Value = 4276
$ws.Range("L97").Value = 4276
$ws.Range("N97").Value = -5268
# Row 100
$ws.Range("H100").Value = 79998.5
$ws.Range("J100").Value = 79998.5
$ws.Range("L100").Value = 79998.5
$ws.Range("N100").Value = -82162.5
# Row 132
$ws.Range("H132").Value = 2175681.2
$ws.Range("I132").Value = 1587.2439
$ws.Range("J132").Value = 20003252
$ws.Range("K132").Value = 4761.7317
$ws.Range("L132").Value = 60009756
$ws.Range("M132").Value = -2231.7317
$ws.Range("N132").Value = -60014816
# Row 136
$ws.Range("H136").Value = 16141088
$ws.Range("I136").Value = 17953138
$ws.Range("K136").Value = 53859414
$ws.Range("M136").Value = -53856864

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5291.467
$ws.Range("I86").Value = 3610.9092
$ws.Range("J86").Value = 9913
$ws.Range("K86").Value = 3610.9092
$ws.Range("L86").Value = 9913
$ws.Range("M86").Value = -2487.9092
$ws.Range("N86").Value = -12159
# Row 89
$ws.Range("H89").Value = 5291.467
$ws.Range("I89").Value = 3610.9092
$ws.Range("J89").Value = 9913
$ws.Range("K89").Value = 18054.546
$ws.Range("L89").Value = 49565
$ws.Range("M89").Value = -12438.546
$ws.Range("N89").Value = -60797
# Row 94
$ws.Range("H94").Value = 1848.1666
$ws.Range("J94").Value = 1365.1111
$ws.Range("L94").Value = 1365.1111
$ws.Range("N94").Value = -2267.1111

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 50065
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
# Row 51
$ws.Range("H51").Value = 19095
$ws.Range("I51").Value = 19095
$ws.Range("K51").Value = 19095
$ws.Range("M51").Value = -18359
# Row 61
$ws.Range("H61").Value = 19095
$ws.Range("I61").Value = 19095
$ws.Range("K61").Value = 19095
$ws.Range("M61").Value = -18747
# Row 86
$ws.Range("H86").Value = 7617.6665
$ws.Range("J86").Value = 9598.143
$ws.Range("L86").Value = 9598.143
$ws.Range("N86").Value = -11844.143
# Row 89
$ws.Range("H89").Value = 7617.6665
$ws.Range("J89").Value = 9598.143
$ws.Range("L89").Value = 47990.715
$ws.Range("N89").Value = -59222.715
# Row 118
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
# Row 132
$ws.Range("H132").Value = 1664.8334
$ws.Range("I132").Value = 1838.2273
$ws.Range("J132").Value = 1188
$ws.Range("K132").Value = 5514.6819
$ws.Range("L132").Value = 3564
$ws.Range("M132").Value = -2984.6819
$ws.Range("N132").Value = -8624

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 13466.2
$ws.Range("J44").Value = 14332.75
$ws.Range("L44").Value = 42998.25
$ws.Range("N44").Value = -43794.25
# Row 107
$ws.Range("H107").Value = 5693592
$ws.Range("J107").Value = 7590189.5
$ws.Range("L107").Value = 22770568.5
$ws.Range("N107").Value = -22774408.5
# Row 114
$ws.Range("H114").Value = 3527.2
$ws.Range("I114").Value = 353.6
$ws.Range("K114").Value = 1060.8
$ws.Range("M114").Value = 2193.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 12737
$ws.Range("I3").Value = 950
$ws.Range("J3").Value = 16666
$ws.Range("K3").Value = 950
$ws.Range("L3").Value = 16666
$ws.Range("M3").Value = -834
$ws.Range("N3").Value = -16898
# Row 132
$ws.Range("H132").Value = 3294914.8
$ws.Range("I132").Value = 3587.3777
$ws.Range("J132").Value = 14687971
$ws.Range("K132").Value = 10762.1331
$ws.Range("L132").Value = 44063913
$ws.Range("M132").Value = -8232.133099999999
$ws.Range("N132").Value = -44068973

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6739.7837
$ws.Range("I7").Value = 6457.0835
$ws.Range("K7").Value = 6457.0835
$ws.Range("M7").Value = -6345.0835
# Row 22
$ws.Range("H22").Value = 10784
$ws.Range("J22").Value = 3089.6365
$ws.Range("L22").Value = 3089.6365
$ws.Range("N22").Value = -3679.6365
# Row 27
$ws.Range("H27").Value = 10784
$ws.Range("J27").Value = 3089.6365
$ws.Range("L27").Value = 3089.6365
$ws.Range("N27").Value = -3303.6365
# Row 40
$ws.Range("H40").Value = 4962.6
$ws.Range("I40").Value = 3786.6
$ws.Range("J40").Value = 8490.6
$ws.Range("K40").Value = 3786.6
$ws.Range("L40").Value = 8490.6
$ws.Range("M40").Value = -3650.6
$ws.Range("N40").Value = -8762.6
# Row 126
$ws.Range("H126").Value = 6739.7837
$ws.Range("I126").Value = 6457.0835
$ws.Range("K126").Value = 19371.2505
$ws.Range("M126").Value = -16901.2505
# Row 132
$ws.Range("H132").Value = 3623.5
$ws.Range("I132").Value = 2078.3794
$ws.Range("K132").Value = 6235.138199999999
$ws.Range("M132").Value = -3705.138199999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 6821.5
$ws.Range("J96").Value = 5298.4287
$ws.Range("L96").Value = 5298.4287
$ws.Range("N96").Value = -8044.4287
# Row 113
$ws.Range("H113").Value = 507.55554
$ws.Range("I113").Value = 257.63635
$ws.Range("J113").Value = 900.2857
$ws.Range("K113").Value = 772.90905
$ws.Range("L113").Value = 2700.8571
$ws.Range("M113").Value = 1397.09095
$ws.Range("N113").Value = -7040.8571
